{"js": "// Update the two-digit multiplication answer table: replace each old\n// \"A\u00d7B=C\" expression with its new value. The mapping below mirrors the\n// unified diff exactly (old text -> new text), applied as independent\n// search-and-replace operations against the document body.\nconst replacements = [\n  [\"96\u00d774=7104\", \"98\u00d759=5782\"],\n  [\"88\u00d727=2376\", \"52\u00d713=676\"],\n  [\"46\u00d718=828\", \"51\u00d763=3213\"],\n  [\"91\u00d777=7007\", \"43\u00d725=1075\"],\n  [\"59\u00d784=4956\", \"82\u00d732=2624\"],\n  [\"95\u00d771=6745\", \"30\u00d758=1740\"],\n  [\"60\u00d759=3540\", \"39\u00d752=2028\"],\n  [\"76\u00d758=4408\", \"95\u00d718=1710\"],\n  [\"12\u00d743=516\", \"31\u00d730=930\"],\n  [\"19\u00d797=1843\", \"14\u00d716=224\"],\n  [\"21\u00d719=399\", \"69\u00d728=1932\"],\n  [\"95\u00d751=4845\", \"94\u00d719=1786\"],\n  [\"84\u00d714=1176\", \"64\u00d716=1024\"],\n  [\"27\u00d758=1566\", \"40\u00d788=3520\"],\n  [\"76\u00d782=6232\", \"23\u00d713=299\"],\n  [\"92\u00d796=8832\", \"97\u00d722=2134\"],\n  [\"18\u00d713=234\", \"40\u00d769=2760\"],\n  [\"83\u00d785=7055\", \"57\u00d725=1425\"],\n  [\"34\u00d750=1700\", \"70\u00d759=4130\"],\n  [\"16\u00d766=1056\", \"33\u00d773=2409\"],\n  [\"75\u00d729=2175\", \"51\u00d784=4284\"],\n  [\"71\u00d797=6887\", \"88\u00d726=2288\"],\n  [\"14\u00d787=1218\", \"56\u00d784=4704\"],\n  [\"50\u00d750=2500\", \"47\u00d796=4512\"],\n  [\"42\u00d766=2772\", \"89\u00d730=2670\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication answer table: replace each old\n# \"A\u00d7B=C\" expression with its new value via Find/Replace across the whole\n# document body (the table cells holding these expressions).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"96\u00d774=7104\"; Replace = \"98\u00d759=5782\" },\n    @{ Find = \"88\u00d727=2376\"; Replace = \"52\u00d713=676\" },\n    @{ Find = \"46\u00d718=828\"; Replace = \"51\u00d763=3213\" },\n    @{ Find = \"91\u00d777=7007\"; Replace = \"43\u00d725=1075\" },\n    @{ Find = \"59\u00d784=4956\"; Replace = \"82\u00d732=2624\" },\n    @{ Find = \"95\u00d771=6745\"; Replace = \"30\u00d758=1740\" },\n    @{ Find = \"60\u00d759=3540\"; Replace = \"39\u00d752=2028\" },\n    @{ Find = \"76\u00d758=4408\"; Replace = \"95\u00d718=1710\" },\n    @{ Find = \"12\u00d743=516\"; Replace = \"31\u00d730=930\" },\n    @{ Find = \"19\u00d797=1843\"; Replace = \"14\u00d716=224\" },\n    @{ Find = \"21\u00d719=399\"; Replace = \"69\u00d728=1932\" },\n    @{ Find = \"95\u00d751=4845\"; Replace = \"94\u00d719=1786\" },\n    @{ Find = \"84\u00d714=1176\"; Replace = \"64\u00d716=1024\" },\n    @{ Find = \"27\u00d758=1566\"; Replace = \"40\u00d788=3520\" },\n    @{ Find = \"76\u00d782=6232\"; Replace = \"23\u00d713=299\" },\n    @{ Find = \"92\u00d796=8832\"; Replace = \"97\u00d722=2134\" },\n    @{ Find = \"18\u00d713=234\"; Replace = \"40\u00d769=2760\" },\n    @{ Find = \"83\u00d785=7055\"; Replace = \"57\u00d725=1425\" },\n    @{ Find = \"34\u00d750=1700\"; Replace = \"70\u00d759=4130\" },\n    @{ Find = \"16\u00d766=1056\"; Replace = \"33\u00d773=2409\" },\n    @{ Find = \"75\u00d729=2175\"; Replace = \"51\u00d784=4284\" },\n    @{ Find = \"71\u00d797=6887\"; Replace = \"88\u00d726=2288\" },\n    @{ Find = \"14\u00d787=1218\"; Replace = \"56\u00d784=4704\" },\n    @{ Find = \"50\u00d750=2500\"; Replace = \"47\u00d796=4512\" },\n    @{ Find = \"42\u00d766=2772\"; Replace = \"89\u00d730=2670\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2) | Out-Null\n}\n"}
